# Cross browser testing: add a "browser" column (C) to the DATA sheet,
# populate it with firefox/chrome values, fix up the execute flags,
# and re-point the hyperlinks that used to live in column C (now D).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# 1) Insert a new blank column at C - this shifts old C/D/E -> D/E/F
#    (cell values, number formats and column widths all move with it).
$ws.Columns.Item(3).Insert()

# 2) Header for the new column
$ws.Range("C1").Value = "browser"

# 3) Fill in the "browser" column + fix the "execute" column (B) to match
#    the new cross-browser test matrix.
$ws.Range("B2").Value = "yes"
$ws.Range("C2").Value = "firefox"

$ws.Range("B3").Value = "no"
$ws.Range("C3").Value = "chrome"

$ws.Range("B4").Value = "no"
$ws.Range("C4").Value = "chrome"

$ws.Range("B5").Value = "no"
$ws.Range("C5").Value = "chrome"

$ws.Range("B6").Value = "yes"
$ws.Range("C6").Value = "chrome"

# 4) Fix the mangled password value that used to live in D6.
$ws.Range("E6").Value = "Ananya@123"

# 5) Rebuild the hyperlinks against their new column (D) locations.
#    (Range.Hyperlinks.Delete() on this engine clears the whole sheet's
#    collection rather than a single cell, so clear once and re-add all
#    five against their new cells; then re-apply the Hyperlink cell style
#    so the cells reuse the original style record instead of a fresh one.)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:ananya111@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:ananya111@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:ananya111@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:trisha@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:trisha@gmail.com")

$ws.Range("D2").Style = "Hyperlink"
$ws.Range("D3").Style = "Hyperlink"
$ws.Range("D4").Style = "Hyperlink"
$ws.Range("D5").Style = "Hyperlink"
$ws.Range("D6").Style = "Hyperlink"

# 6) Update the sheet selection to match the new layout.
$ws.Range("C3:C5").Select()
